$wb = $excel.ActiveWorkbook

# "Data quality - scores" sheet: the "Device completeness" column (B) previously
# duplicated the retention sentence that already lives in the "Retention" column
# (F). Trim column B down to just the device-availability statement.
$ws = $wb.Worksheets.Item("Data quality - scores")

$ws.Range("B4").Value = "Data source is available from `n0-25% of the devices."
$ws.Range("B5").Value = "Data source is available from 26-50% of the devices."
$ws.Range("B6").Value = "Data source is available from 51-75% of the devices."
$ws.Range("B7").Value = "Data source is available from 76-100% of the devices."
$ws.Range("B8").Value = "Data source is available for 100% of the devices.`n"

# The "Data quality - scores" tab is now the active/selected tab, with H4 as the
# active cell (previously "Detection scores - descriptions" was active).
$ws.Activate()
$ws.Range("H4").Select()
